$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "최종점수" (K column) values
$ws.Range("K2").Value = 59.7
$ws.Range("K3").Value = 57.3
$ws.Range("K4").Value = 50.5
$ws.Range("K5").Value = 47.7

# Update "MACRO_SCORE" (N column) values
$ws.Range("N2").Value = 51.15965480231979
$ws.Range("N3").Value = 51.15965480231979
$ws.Range("N4").Value = 51.15965480231979
$ws.Range("N5").Value = 51.15965480231979
